# Update stack-trace line numbers in the sample generated output text,
# reflecting upstream source line shifts caused by:
# "Fixed #295 Add the version of M2Doc in the template custom properties."

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# caseLet(M2DocEvaluator.java:1050) -> 1096
Replace-Text "caseLet(M2DocEvaluator.java:1050)" "caseLet(M2DocEvaluator.java:1096)"

# doSwitch(M2DocEvaluator.java:1038) -> 1084 (appears 3 times)
Replace-Text "doSwitch(M2DocEvaluator.java:1038)" "doSwitch(M2DocEvaluator.java:1084)"
Replace-Text "doSwitch(M2DocEvaluator.java:1038)" "doSwitch(M2DocEvaluator.java:1084)"
Replace-Text "doSwitch(M2DocEvaluator.java:1038)" "doSwitch(M2DocEvaluator.java:1084)"

# caseBlock(M2DocEvaluator.java:1254) -> 1300
Replace-Text "caseBlock(M2DocEvaluator.java:1254)" "caseBlock(M2DocEvaluator.java:1300)"

# caseDocumentTemplate(M2DocEvaluator.java:275) -> 278
Replace-Text "caseDocumentTemplate(M2DocEvaluator.java:275)" "caseDocumentTemplate(M2DocEvaluator.java:278)"

# generate(M2DocEvaluator.java:264) -> 267
Replace-Text "generate(M2DocEvaluator.java:264)" "generate(M2DocEvaluator.java:267)"

# generate(M2DocUtils.java:712) -> 694
Replace-Text "generate(M2DocUtils.java:712)" "generate(M2DocUtils.java:694)"

# prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459) -> 475
Replace-Text "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)" "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)"

# generation(AbstractTemplatesTestSuite.java:369) -> 384
Replace-Text "generation(AbstractTemplatesTestSuite.java:369)" "generation(AbstractTemplatesTestSuite.java:384)"
